# Add a new "AVPC" (Annual Vehicle Parking Cost) row to the "Key to Variables"
# sheet, right after the existing "AVMC" (Annual Vehicle Maintenance Cost) row.
# Commit message: "Add vehicle parking costs (#161)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# The new row goes at sheet row 191 (just below the existing AVMC row, 190),
# pushing every row currently at 191..232 down by one (to 192..233).
$newRow = 191

# Insert a blank row; formatting is pulled down from the row above (190),
# which already has the right "Top Level Folder" (A) styling for this block.
$ws.Rows.Item($newRow).Insert()

# Fill in the new row's values.
$ws.Cells.Item($newRow, 1).Value = "trans"
$ws.Cells.Item($newRow, 2).Value = "AVPC"
$ws.Cells.Item($newRow, 3).Value = "Annual Vehicle Parking Cost"
$ws.Cells.Item($newRow, 6).Value = "medium"

# The "Importance to Update" cell (column F) needs the "medium" fill used by
# other medium-importance rows, rather than the "low" fill it inherited from
# row 190 above. Copy that single cell's formatting from a known medium row
# (row 199, after the insert shifted the original row 198 down by one).
$ws.Range("F199").Copy()
$ws.Range("F191").PasteSpecial(-4122)
$excel.CutCopyMode = 0
